# Bug fix: headers used underscores (Start_Date, End_Date, Pending_days)
# which made the summary look like it was repeating the same bank row;
# rename them to the human readable "Start Date", "End Date", "Pending days".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Start Date"
$ws.Range("D1").Value = "End Date"
$ws.Range("E1").Value = "Pending days"

# Move the active selection to E1 (matches the resulting workbook view state)
$ws.Range("E1").Select()
